# Apply updates described by the commit "Add files via upload"
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Лист1")
$ws2 = $wb.Worksheets.Item("GirHun")

# Insert two new rows into sheet1 (Лист1): new row at position 10 and new row at position 16 (post-insert numbering)
$ws1.Rows.Item(10).Insert()
$ws1.Rows.Item(16).Insert()

# Row 10 (new): rank 52374, user NeaI_Wu, Python=605, Python3=61, MySQL=7, Github=No data
$ws1.Range("A10").Value = 52374
$ws1.Range("A10").NumberFormat = "#,##0"
$ws1.Range("B10").Value = "https://leetcode.com/u/NeaI_Wu/"
$ws1.Range("C10").Value = 605
$ws1.Range("D10").Value = 61
$ws1.Range("F10").Value = 7
$ws1.Range("L10").Value = "No data"

# Row 16 (new): rank 143368, user tapanvaishnav17, Python=22, Python3=351, MySQL=1, Github=No data
$ws1.Range("A16").Value = 143368
$ws1.Range("A16").NumberFormat = "#,##0"
$ws1.Range("B16").Value = "https://leetcode.com/u/tapanvaishnav17/"
$ws1.Range("C16").Value = 22
$ws1.Range("D16").Value = 351
$ws1.Range("F16").Value = 1
$ws1.Range("L16").Value = "No data"

# Update the active sheet / selection bookkeeping
$ws1.Activate()
$ws1.Range("L17").Select()
